# Generate Report for Archive
#
# The localization status report's "Status" column used to read
# "Ready for handoff" for every row; regenerating the report now captures
# these files mid-flight, so the status text becomes "In Translation".
# That shared text is shown on all three tabs:
#   - Overview!E2:F4  (the zh-cn / de-de roll-up "status" columns)
#   - zh-cn!C2:C4     (the "Status" column of the zh-cn detail table)
#   - de-de!C2:C4     (the "Status" column of the de-de detail table)
#
# Because the new status string is shorter than the old one, the "Status"
# columns on all three sheets shrink to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update every cell that currently shows the old "Ready for handoff" status
# text to the new status text.
$wsOverview.Range("E2:F4").Value = $newStatus
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsDeDe.Range("C2:C4").Value = $newStatus

# Re-fit the "Status" columns now that the text is shorter.
$fitWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $fitWidth
$wsOverview.Columns.Item(6).ColumnWidth = $fitWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $fitWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $fitWidth
